$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column text values that look numeric are written as text, matching
# the original inline-string cell type, by temporarily forcing a text number
# format on the D column, then clearing the format afterwards so the cells
# end up with no explicit style (same as before the edit).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '58.983.71'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').Value = '2.586.75'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '520.91'
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('D6').Value = '139.76'
$ws.Range('E6').Value = '  -2.42%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.44%  '
$ws.Range('D9').Value = '2.597.78'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').Value = '6.55'
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('E12').Value = '  +1.50%  '
$ws.Range('E13').Value = '  +3.12%  '
$ws.Range('D14').Value = '3.039.40'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').Value = '58.944.13'
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').Value = '2.575.62'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').Value = '338.56'
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').Value = '4.30'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('E21').Value = '  -1.32%  '
$ws.Range('D22').Value = '6.52'
$ws.Range('E22').Value = '  +2.80%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = '66.13'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '0.0₃0726'
$ws.Range('E30').Value = '  -3.74%  '
$ws.Range('D31').Value = '5.95'
$ws.Range('E31').Value = '  -6.36%  '
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').Value = '18.70'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '148.97'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('E36').Value = '  -2.36%  '
$ws.Range('D37').Value = '36.81'
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('D38').Value = '1.47'
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').Value = '0.811'
$ws.Range('E40').Value = '  -7.13%  '
$ws.Range('D41').Value = '3.51'
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = '272.69'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('D44').Value = '10.77'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('D46').Value = '0.0950'
$ws.Range('E46').Value = '  -0.50%  '
$ws.Range('D47').Value = '0.0518'
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('D48').Value = '18.47'
$ws.Range('E48').Value = '  -2.09%  '
$ws.Range('D49').Value = '1.972.63'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('E51').Value = '  -0.51%  '

$ws.Range("D2:D51").ClearFormats()
